# Reorders the "Steps" content of test cases TC1, TC3 and TC4 so that:
#   TC1 (rows 10-11)  <- old TC3 content (Periodos Avaliativos)
#   TC3 (rows 28-29)  <- old TC4 content (Avaliacoes)
#   TC4 (rows 37-38)  <- old TC1 content (Competencias (portfolio))
# TC2 (rows 19-20) and TC5 (rows 46-47) stay unchanged.
# This corresponds to version bump 1.0 -> 1.1 described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original (1.0) text values before overwriting anything.
# NOTE: use Value2 for reads (Value getter is unreliable in this COM shim).
$tc1_b10 = $ws.Range("B10").Value2
$tc1_d10 = $ws.Range("D10").Value2
$tc1_b11 = $ws.Range("B11").Value2

$tc3_b28 = $ws.Range("B28").Value2
$tc3_d28 = $ws.Range("D28").Value2
$tc3_b29 = $ws.Range("B29").Value2

$tc4_b37 = $ws.Range("B37").Value2
$tc4_d37 = $ws.Range("D37").Value2
$tc4_b38 = $ws.Range("B38").Value2

# TC1 slot now gets the old TC3 (Periodos Avaliativos) content
$ws.Range("B10").Value = $tc3_b28
$ws.Range("D10").Value = $tc3_d28
$ws.Range("B11").Value = $tc3_b29

# TC3 slot now gets the old TC4 (Avaliacoes) content
$ws.Range("B28").Value = $tc4_b37
$ws.Range("D28").Value = $tc4_d37
$ws.Range("B29").Value = $tc4_b38

# TC4 slot now gets the old TC1 (Competencias (portfolio)) content
$ws.Range("B37").Value = $tc1_b10
$ws.Range("D37").Value = $tc1_d10
$ws.Range("B38").Value = $tc1_b11
